# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# "展览" sheet (first sheet): update "want to go" counts (column F)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 111
$wsExhibit.Range("F5").Value = 2950
$wsExhibit.Range("F6").Value = 298
$wsExhibit.Range("F7").Value = 399

# "全部类型" sheet (aggregated list, fourth sheet): same rows, different
# row positions because it merges multiple categories
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 111
$wsAll.Range("F5").Value = 2950
$wsAll.Range("F6").Value = 298
$wsAll.Range("F9").Value = 399
